$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '29.168.85'
$c.ClearFormats()
$ws.Range("E2").Value = '  -0.56%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.824.14'
$c.ClearFormats()
$ws.Range("E3").Value = '  -0.85%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c.ClearFormats()
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '235.68'
$c.ClearFormats()
$ws.Range("E5").Value = '  -1.46%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.6095'
$c.ClearFormats()
$ws.Range("E6").Value = '  -2.93%  '
$ws.Range("E7").Value = '  +0.12%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.07091'
$c.ClearFormats()
$ws.Range("E8").Value = '  -4.45%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.2802'
$c.ClearFormats()
$ws.Range("E9").Value = '  -2.98%  '
$ws.Range("E10").Value = '  -5.86%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07660'
$c.ClearFormats()
$ws.Range("E11").Value = '  -0.89%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.823.65'
$c.ClearFormats()
$ws.Range("E12").Value = '  +1.16%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.798'
$c.ClearFormats()
$ws.Range("E13").Value = '  -3.09%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.000009999'
$c.ClearFormats()
$ws.Range("E14").Value = '  -1.56%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.6319'
$c.ClearFormats()
$ws.Range("E15").Value = '  -6.16%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.069.34'
$c.ClearFormats()
$ws.Range("E16").Value = '  -0.83%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '78.58'
$c.ClearFormats()
$ws.Range("E17").Value = '  -3.61%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '5.864'
$c.ClearFormats()
$ws.Range("E18").Value = '  -5.51%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '29.168.04'
$c.ClearFormats()
$ws.Range("E19").Value = '  -0.53%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '226.56'
$c.ClearFormats()
$ws.Range("E20").Value = '  -0.87%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  -4.24%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.989'
$c.ClearFormats()
$ws.Range("E23").Value = '  -4.57%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.9997'
$c.ClearFormats()
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("E25").Value = '  -1.45%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.034'
$c.ClearFormats()
$ws.Range("E26").Value = '  -4.96%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.1307'
$c.ClearFormats()
$ws.Range("E27").Value = '  -2.79%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '16.57'
$c.ClearFormats()
$ws.Range("E28").Value = '  -4.43%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.485'
$c.ClearFormats()
$ws.Range("E29").Value = '  +1.80%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.06297'
$c.ClearFormats()
$ws.Range("E30").Value = '  -14.48%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.450'
$c.ClearFormats()
$ws.Range("E31").Value = '  -1.54%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.825'
$c.ClearFormats()
$ws.Range("E32").Value = '  -4.97%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.803'
$c.ClearFormats()
$ws.Range("E33").Value = '  -5.74%  '
$ws.Range("E34").Value = '  -1.16%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.737'
$c.ClearFormats()
$ws.Range("E35").Value = '  -4.47%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.6427'
$c.ClearFormats()
$ws.Range("E36").Value = '  -7.22%  '
$ws.Range("E37").Value = '  -1.26%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.213.47'
$c.ClearFormats()
$ws.Range("E38").Value = '  -1.45%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.720'
$c.ClearFormats()
$ws.Range("E39").Value = '  -3.01%  '
$ws.Range("E40").Value = '  -5.21%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.533'
$c.ClearFormats()
$ws.Range("E41").Value = '  -5.33%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.9067'
$c.ClearFormats()
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("E43").Value = '  +0.12%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '100.76'
$c.ClearFormats()
$ws.Range("E44").Value = '  +0.16%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.976.90'
$c.ClearFormats()
$ws.Range("E45").Value = '  +0.11%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '62.66'
$c.ClearFormats()
$ws.Range("E46").Value = '  -3.85%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.00000000117'
$c.ClearFormats()
$ws.Range("E47").Value = '  -2.76%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.605'
$c.ClearFormats()
$ws.Range("E48").Value = '  -5.50%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '8.502'
$c.ClearFormats()
$ws.Range("E49").Value = '  -4.07%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.4560'
$c.ClearFormats()
$ws.Range("E50").Value = '  -0.60%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.05512'
$c.ClearFormats()
$ws.Range("E51").Value = '  -2.68%  '
